# Update event queue diagrams (slide 1)
# Renames the example event-listener / event labels in the
# "EVENT LISTENERS" diagram:
#   Event Listener 1 -> Event Listener A
#   Event Z           -> Event N
#   Event A           -> Event 1
#   Event Listener N  -> Event Listener Z

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The diagram shapes live inside a (nested) group, but GroupItems
# exposes them by their unique shape name regardless of nesting depth.
$group = $s.Shapes.Item(1)
$items = $group.GroupItems

# 1) "Event Listener 1" -> "Event Listener A"  (single run, simple swap)
$items.Item("Rectangle 39").TextFrame.TextRange.Text = "Event Listener A"

# 2) "Event " + "Z" (two runs) -> single run "Event N"
#    Re-typing the whole box's text first with unrelated placeholder
#    content forces the box back down to one run (keeping the first
#    run's formatting), then we set the final wording.
$eventZBox = $items.Item("ZoneTexte 67")
$eventZBox.TextFrame.TextRange.Text = "zzzzzzzzzzzzzzzzzzzzzzzz"
$eventZBox.TextFrame.TextRange.Text = "Event N"

# 3) "Event A" -> "Event 1"
$items.Item("ZoneTexte 70").TextFrame.TextRange.Text = "Event 1"

# 4) "Event Listener N" -> "Event Listener Z"
$items.Item("Rectangle 50").TextFrame.TextRange.Text = "Event Listener Z"
